$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.24 = 49555.03 pesos`n✅ 49555.03 pesos = 12.24 = 980.59 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- tasas: update rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 81.69
$wsTasas.Range("O10").Value = 4048.15
$wsTasas.Range("N12").Value = 4050
$wsTasas.Range("O12").Value = 80.14100000000001
